$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# New lesson breakdown rows for the "Email" topic (D19:D25).
# Cell text is set in this specific order so that the new shared-string
# table entries come out in the same order as the target workbook.
$ws.Range("D19").Value = "Emailes veszélyek"
$ws.Range("D20").Value = "Email phising gyakori kinézete"
$ws.Range("D22").Value = "Hogyan ismerjük fel ezeket a veszélyeket"
$ws.Range("D21").Value = "Email spam gyakori kinézete"
$ws.Range("D23").Value = "Phising weboldal példa"
$ws.Range("D24").Value = "Hogyan védekezzünk - ""Józan paraszti ész használata"", ""semmi sincs ingyen"""
$ws.Range("D25").Value = "Mire tudják felhasználni az adatainkat, ha megszerzik"

# New detail text for the "Free wifi" topic (E7:E8).
$ws.Range("E7").Value = "Free wifi - dejó hogy van de elmondani miért veszélyes"
$ws.Range("E8").Value = "Példa arra miért veszélyes, hogyan tudnak átverni/meglopni ezzel"

# D21 and D23 pick up a slightly different cell style (explicit "no fill"
# combined with the existing border) - reproduce it on D21 then copy the
# resulting format onto D23 so both cells share the same style record.
$ws.Range("D21").Interior.Color = 65535
$ws.Range("D21").Interior.Pattern = -4142
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column E grew a bit wider to fit the new content.
$ws.Columns.Item(5).ColumnWidth = 60

# Update the view: scrolled so column E is visible and F7 is selected.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F7").Select() | Out-Null
